# Refresh the "cryptos" price list (Price / Volume(1h) columns, plus a
# couple of re-ranked coin rows) to match the latest GitHub Actions scrape.
#
# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (losing e.g. trailing zeros, or turning "4.619" into the number 4.619
# instead of the literal label) are first switched to a Text number format
# so the assignment sticks as the exact original string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.560.23"
$ws.Range("E2").Value = "  +5.45%  "

$ws.Range("D3").Value = "1.724.31"
$ws.Range("E3").Value = "  +4.26%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.02"
$ws.Range("E5").Value = "  +3.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5401"
$ws.Range("E6").Value = "  +3.22%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2683"
$ws.Range("E8").Value = "  +1.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06615"
$ws.Range("E9").Value = "  +4.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.82"
$ws.Range("E10").Value = "  +6.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07732"
$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.619"

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.739.34"
$ws.Range("E13").Value = "  +7.39%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.961.91"
$ws.Range("E14").Value = "  +4.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5871"
$ws.Range("E15").Value = "  +4.60%  "

$ws.Range("D16").Value = "0.0₅8312"
$ws.Range("E16").Value = "  +1.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.04"
$ws.Range("E17").Value = "  +3.91%  "

$ws.Range("D18").Value = "27.568.67"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.34"
$ws.Range("E19").Value = "  +15.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.739"
$ws.Range("E21").Value = "  +1.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.100"
$ws.Range("E23").Value = "  +2.55%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.31"
$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1237"
$ws.Range("E26").Value = "  +3.56%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.696"
$ws.Range("E27").Value = "  +12.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.416"
$ws.Range("E28").Value = "  +2.14%  "

$ws.Range("E29").Value = "  +4.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05557"
$ws.Range("E30").Value = "  +1.86%  "

$ws.Range("E31").Value = "  +2.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.550"
$ws.Range("E32").Value = "  +2.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.460"
$ws.Range("E33").Value = "  +2.59%  "

$ws.Range("E34").Value = "  +6.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9642"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.824"
$ws.Range("E36").Value = "  +1.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.446"
$ws.Range("E37").Value = "  +1.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5960"
$ws.Range("E38").Value = "  +5.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01650"
$ws.Range("E39").Value = "  +4.32%  "

$ws.Range("E40").Value = "  +1.15%  "

$ws.Range("D41").Value = "1.056.45"
$ws.Range("E41").Value = "  +2.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8539"

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.47"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("D45").Value = "1.867.27"
$ws.Range("E45").Value = "  +4.13%  "

$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").Value = "  +9.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.16"
$ws.Range("E47").Value = "  +2.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.194"
$ws.Range("E48").Value = "  +2.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4442"
$ws.Range("E49").Value = "  +2.38%  "

$ws.Range("E50").Value = "  +0.70%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05278"
$ws.Range("E51").Value = "  +1.62%  "
